$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC4").Value = 13
$ws.Range("AE4").Value = 12
$ws.Range("AL4").Value = 19
$ws.Range("G4").Value = 2.6
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 2.63
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 4
$ws.Range("AH5").Value = 13
$ws.Range("AA7").Value = 22
$ws.Range("AB7").Value = 37
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 5.9
$ws.Range("AE7").Value = 16.5
$ws.Range("AI7").Value = 14
$ws.Range("AP7").Value = 22
$ws.Range("AQ7").Value = 55
$ws.Range("AT7").Value = 2.35
$ws.Range("J7").Value = 2.95
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 7.4
$ws.Range("P7").Value = 2.45
$ws.Range("Q7").Value = 2.22
$ws.Range("R7").Value = 1.52
$ws.Range("S7").Value = 1.45
$ws.Range("T7").Value = 2.37
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 1.7
$ws.Range("W7").Value = 6.6
$ws.Range("X7").Value = 10.5
$ws.Range("Z7").Value = 24
$ws.Range("AA9").Value = 8.75
$ws.Range("AH9").Value = 17.5
$ws.Range("AO9").Value = 5.8
$ws.Range("AP9").Value = 14.5
$ws.Range("AQ9").Value = 15
$ws.Range("AR9").Value = 37
$ws.Range("AT9").Value = 3.35
$ws.Range("AY9").Value = 40
$ws.Range("AZ9").Value = 300
$ws.Range("I9").Value = 7.1
$ws.Range("J9").Value = 1.75
$ws.Range("K9").Value = 2.5
$ws.Range("L9").Value = 6.4
$ws.Range("Q9").Value = 1.52
$ws.Range("R9").Value = 2.22
$ws.Range("U9").Value = 1.82
$ws.Range("V9").Value = 1.94
$ws.Range("W9").Value = 6.9
$ws.Range("X9").Value = 6.1
$ws.Range("Z9").Value = 7.5
$ws.Range("AB10").Value = 40
$ws.Range("AC10").Value = 13.5
$ws.Range("AD10").Value = 7.5
$ws.Range("AH10").Value = 6.6
$ws.Range("AJ10").Value = 7
$ws.Range("AK10").Value = 8.25
$ws.Range("AN10").Value = 7.3
$ws.Range("AO10").Value = 32
$ws.Range("AR10").Value = 200
$ws.Range("AT10").Value = 3.1
$ws.Range("AX10").Value = 6.5
$ws.Range("AZ10").Value = 18.5
$ws.Range("BA10").Value = 45
$ws.Range("G10").Value = 5.7
$ws.Range("H10").Value = 4.25
$ws.Range("I10").Value = 1.44
$ws.Range("J10").Value = 5.5
$ws.Range("K10").Value = 2.37
$ws.Range("L10").Value = 1.91
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 11.4
$ws.Range("O10").Value = 1.17
$ws.Range("P10").Value = 4.05
$ws.Range("Q10").Value = 1.6
$ws.Range("R10").Value = 2.07
$ws.Range("S10").Value = 1.3
$ws.Range("T10").Value = 3.32
$ws.Range("U10").Value = 1.75
$ws.Range("V10").Value = 2.03
$ws.Range("W10").Value = 14.5
$ws.Range("X10").Value = 29
$ws.Range("Y10").Value = 15
$ws.Range("Z10").Value = 80
$ws.Range("AB11").Value = 175
$ws.Range("AC11").Value = 17.5
$ws.Range("AD11").Value = 24
$ws.Range("AE11").Value = 70
$ws.Range("AF11").Value = 350
$ws.Range("AH11").Value = 8.75
$ws.Range("AI11").Value = 5.4
$ws.Range("AJ11").Value = 14
$ws.Range("AK11").Value = 4.7
$ws.Range("AL11").Value = 11.75
$ws.Range("AM11").Value = 55
$ws.Range("AN11").Value = 37
$ws.Range("AO11").Value = 500
$ws.Range("AP11").Value = 250
$ws.Range("AU11").Value = 16.5
$ws.Range("AV11").Value = 200
$ws.Range("AW11").Value = 2.95
$ws.Range("AX11").Value = 3.65
$ws.Range("AY11").Value = 19.5
$ws.Range("AZ11").Value = 6.5
$ws.Range("BA11").Value = 35
$ws.Range("BB11").Value = 350
$ws.Range("G11").Value = 37
$ws.Range("I11").Value = 1.04
$ws.Range("J11").Value = 27
$ws.Range("K11").Value = 3.35
$ws.Range("L11").Value = 1.26
$ws.Range("O11").Value = 1.05
$ws.Range("P11").Value = 7.9
$ws.Range("R11").Value = 3.48
$ws.Range("S11").Value = 1.14
$ws.Range("T11").Value = 4.9
$ws.Range("U11").Value = 2.93
$ws.Range("V11").Value = 1.38
$ws.Range("W11").Value = 120
$ws.Range("X11").Value = 200
$ws.Range("Y11").Value = 200
$ws.Range("AJ13").Value = 10
$ws.Range("AK13").Value = 26
$ws.Range("AL13").Value = 21
$ws.Range("AM13").Value = 29
$ws.Range("AO13").Value = 15
$ws.Range("AQ13").Value = 51
$ws.Range("AY13").Value = 23
$ws.Range("G13").Value = 2.5
$ws.Range("I13").Value = 2.7
$ws.Range("L13").Value = 3.4
$ws.Range("AB14").Value = 34
$ws.Range("AF14").Value = 81
$ws.Range("AK14").Value = 101
$ws.Range("AN14").Value = 3.2
$ws.Range("AR14").Value = 51
$ws.Range("AT14").Value = 2.75
$ws.Range("AU14").Value = 10
$ws.Range("AZ14").Value = 201
$ws.Range("G14").Value = 1.4
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = 8
$ws.Range("L14").Value = 7.5
$ws.Range("Q14").Value = 1.93
$ws.Range("R14").Value = 1.93
$ws.Range("S14").Value = 1.4
$ws.Range("T14").Value = 2.75
$ws.Range("U14").Value = 2.2
$ws.Range("V14").Value = 1.62
$ws.Range("X14").Value = 6
$ws.Range("Y14").Value = 9
$ws.Range("AI16").Value = 12
$ws.Range("AO16").Value = 17
$ws.Range("G16").Value = 2.6
$ws.Range("I16").Value = 2.8
$ws.Range("J16").Value = 3.5
$ws.Range("L16").Value = 3.6
$ws.Range("U16").Value = 2
$ws.Range("V16").Value = 1.73
$ws.Range("W16").Value = 7
$ws.Range("AA17").Value = 12.5
$ws.Range("AB17").Value = 16
$ws.Range("AD17").Value = 9.25
$ws.Range("AE17").Value = 11.5
$ws.Range("AF17").Value = 29
$ws.Range("AJ17").Value = 12.5
$ws.Range("AK17").Value = 50
$ws.Range("AM17").Value = 22
$ws.Range("AP17").Value = 12.5
$ws.Range("AQ17").Value = 25
$ws.Range("AX17").Value = 17
$ws.Range("BB17").Value = 120
$ws.Range("L17").Value = 3.65
$ws.Range("P17").Value = 5.5
$ws.Range("Q17").Value = 1.36
$ws.Range("R17").Value = 2.9
$ws.Range("U17").Value = 1.38
$ws.Range("V17").Value = 2.8
$ws.Range("W17").Value = 14
$ws.Range("X17").Value = 13.5
$ws.Range("Z17").Value = 18.5
